$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '29.076.20'
Set-TextValue 'E2' '  -0.42%  '
Set-TextValue 'D3' '1.815.01'
Set-TextValue 'E3' '  -0.79%  '
Set-TextValue 'E5' '  -1.95%  '
Set-TextValue 'D6' '0.5851'
Set-TextValue 'E6' '  -3.45%  '
Set-TextValue 'E7' '  +0.64%  '
Set-TextValue 'D8' '0.2706'
Set-TextValue 'E8' '  -4.26%  '
Set-TextValue 'D9' '0.06732'
Set-TextValue 'E9' '  -5.09%  '
Set-TextValue 'D10' '22.68'
Set-TextValue 'E10' '  -5.50%  '
Set-TextValue 'D11' '0.07537'
Set-TextValue 'E11' '  -1.35%  '
Set-TextValue 'D12' '1.813.45'
Set-TextValue 'E12' '  -0.87%  '
Set-TextValue 'D13' '4.620'
Set-TextValue 'E13' '  -3.75%  '
Set-TextValue 'D14' '0.6131'
Set-TextValue 'E14' '  -3.80%  '
Set-TextValue 'D15' '0.000009337'
Set-TextValue 'E15' '  -5.97%  '
Set-TextValue 'D16' '74.07'
Set-TextValue 'E16' '  -7.16%  '
Set-TextValue 'D17' '28.841.99'
Set-TextValue 'E17' '  -1.13%  '
Set-TextValue 'D18' '5.388'
Set-TextValue 'E18' '  -10.02%  '
Set-TextValue 'E19' '  +0.54%  '
Set-TextValue 'D20' '206.13'
Set-TextValue 'E20' '  -10.60%  '
Set-TextValue 'D21' '11.31'
Set-TextValue 'E21' '  -4.16%  '
Set-TextValue 'D22' '6.709'
Set-TextValue 'E22' '  -4.15%  '
Set-TextValue 'D23' '1.008'
Set-TextValue 'E23' '  +0.70%  '
Set-TextValue 'D24' '154.17'
Set-TextValue 'E24' '  -0.84%  '
Set-TextValue 'D25' '7.732'
Set-TextValue 'E25' '  -3.90%  '
Set-TextValue 'D26' '0.1252'
Set-TextValue 'E26' '  -2.94%  '
Set-TextValue 'D27' '16.12'
Set-TextValue 'E27' '  -3.50%  '
Set-TextValue 'D28' '1.408'
Set-TextValue 'E28' '  -3.05%  '
Set-TextValue 'D29' '0.06292'
Set-TextValue 'E29' '  -5.66%  '
Set-TextValue 'D30' '1.429'
Set-TextValue 'E30' '  -2.05%  '
Set-TextValue 'D31' '3.680'
Set-TextValue 'E31' '  -3.49%  '
Set-TextValue 'D32' '3.656'
Set-TextValue 'E32' '  -4.61%  '
Set-TextValue 'E33' '  -2.81%  '
Set-TextValue 'D34' '1.039'
Set-TextValue 'E34' '  -8.10%  '
Set-TextValue 'D35' '2.532'
Set-TextValue 'E35' '  -0.59%  '
Set-TextValue 'D36' '0.6294'
Set-TextValue 'E36' '  -4.19%  '
Set-TextValue 'D37' '2.749'
Set-TextValue 'E37' '  -0.24%  '
Set-TextValue 'D38' '0.01707'
Set-TextValue 'E38' '  -3.27%  '
Set-TextValue 'D39' '6.363'
Set-TextValue 'E39' '  -3.26%  '
Set-TextValue 'D40' '1.126.06'
Set-TextValue 'E40' '  -8.97%  '
Set-TextValue 'D41' '0.8602'
Set-TextValue 'E41' '  -7.18%  '
Set-TextValue 'E42' '  +0.60%  '
Set-TextValue 'D43' '1.963.97'
Set-TextValue 'E43' '  -0.96%  '
Set-TextValue 'D44' '99.73'
Set-TextValue 'E44' '  -0.52%  '
Set-TextValue 'B45' 'Aave'
Set-TextValue 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '59.79'
Set-TextValue 'E45' '  -5.84%  '
Set-TextValue 'B46' 'BabyDogeCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.00000000113'
Set-TextValue 'E46' '  -2.88%  '
Set-TextValue 'D47' '0.4544'
Set-TextValue 'E47' '  -0.41%  '
Set-TextValue 'B48' 'Cronos'
Set-TextValue 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.05490'
Set-TextValue 'E48' '  -1.69%  '
Set-TextValue 'B49' 'RenderToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D49' '1.552'
Set-TextValue 'E49' '  -4.75%  '
Set-TextValue 'D50' '8.188'
Set-TextValue 'E50' '  -3.80%  '
Set-TextValue 'D51' '0.9971'
Set-TextValue 'E51' '  -0.01%  '
